# Auto-generated Excel COM-interop script to apply the bluenile.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "Hyperlink" cell style (currently style index 2, used by B2)
# by copying its formatting to a scratch cell far outside the used range, so we can
# reapply the exact same style later without Excel fabricating a brand-new style entry.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null

# Remove the old hyperlinks and wipe all existing data (values + cell formatting)
# in the table area so we can rebuild it cleanly to match the new layout.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Range("A1:E28").Clear() | Out-Null

# --- Header row ---
$ws.Range("A1").Value = "action"
$ws.Range("B1").Value = "locator"
$ws.Range("C1").Value = "value"
$ws.Range("D1").Value = "waitBefore"
$ws.Range("E1").Value = "waitAfter"

# --- Data rows 2-27 ---
# Row 2
$ws.Range("A2").Value = "goto"
$ws.Range("B2").Value = "https://www.bluenile.com/jewelry/necklaces/lab-grown-diamond-cushion-cut-solitaire-pendant-in-14k-white-gold-1-2-ct-tw-f-g-vs2-si1-item-202314"
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 9000

# Row 3
$ws.Range("A3").Value = "scroll"
$ws.Range("B3").Value = "Ships by"
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 5000

# Row 4
$ws.Range("A4").Value = "click"
$ws.Range("B4").Value = "ADD TO CART button"
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 2000

# Row 5
$ws.Range("A5").Value = "goto"
$ws.Range("B5").Value = "https://www.bluenile.com/shopping-cart"
$ws.Range("D5").Value = 1000
$ws.Range("E5").Value = 9000

# Row 6
$ws.Range("A6").Value = "wairfortext"
$ws.Range("B6").Value = "Summary"

# Row 7
$ws.Range("A7").Value = "scroll"
$ws.Range("B7").Value = "Book Now button on popup"
$ws.Range("D7").Value = 1000
$ws.Range("E7").Value = 5000

# Row 8
$ws.Range("A8").Value = "click"
$ws.Range("B8").Value = "Book Now button on popup"
$ws.Range("D8").Value = 1000
$ws.Range("E8").Value = 5000

# Row 9
$ws.Range("A9").Value = "scroll"
$ws.Range("B9").Value = "We Accept"
$ws.Range("D9").Value = 1000
$ws.Range("E9").Value = 5000
$ws.Range("B9").Style = "Normal"

# Row 10
$ws.Range("A10").Value = "click"
$ws.Range("B10").Value = "Checkout button"
$ws.Range("D10").Value = 1000
$ws.Range("E10").Value = 5000

# Row 11
$ws.Range("A11").Value = "wairfortext"
$ws.Range("B11").Value = "Please provide an email address"

# Row 12
$ws.Range("A12").Value = "click"
$ws.Range("B12").Value = "Email Address input field"
$ws.Range("D12").Value = 1000
$ws.Range("E12").Value = 2000

# Row 13
$ws.Range("A13").Value = "type"
$ws.Range("B13").Value = "Email Address input field"
$ws.Range("C13").Value = "mellina@gmail.com"
$ws.Range("D13").Value = 1000
$ws.Range("E13").Value = 2000

# Row 14
$ws.Range("A14").Value = "click"
$ws.Range("B14").Value = "Continue button"
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = 5000

# Row 15
$ws.Range("A15").Value = "wairfortext"
$ws.Range("B15").Value = "First Name"

# Row 16
$ws.Range("A16").Value = "click"
$ws.Range("B16").Value = "First Name input field"
$ws.Range("D16").Value = 1000
$ws.Range("E16").Value = 2000

# Row 17
$ws.Range("A17").Value = "type"
$ws.Range("B17").Value = "First Name input field"
$ws.Range("C17").Value = "Mellina"
$ws.Range("D17").Value = 1000
$ws.Range("E17").Value = 2000

# Row 18
$ws.Range("A18").Value = "click"
$ws.Range("B18").Value = "Last Name input field"
$ws.Range("D18").Value = 1000
$ws.Range("E18").Value = 2000

# Row 19
$ws.Range("A19").Value = "type"
$ws.Range("B19").Value = "Last Name input field"
$ws.Range("C19").Value = "James"
$ws.Range("D19").Value = 1000
$ws.Range("E19").Value = 2000

# Row 20
$ws.Range("A20").Value = "click"
$ws.Range("B20").Value = "Enter Address Manually button"
$ws.Range("D20").Value = 1000
$ws.Range("E20").Value = 2000

# Row 21
$ws.Range("A21").Value = "scroll"
$ws.Range("B21").Value = "Zip / Postal Code"
$ws.Range("D21").Value = 1000
$ws.Range("E21").Value = 3000

# Row 22
$ws.Range("A22").Value = "filldata"
$ws.Range("B22").Value = "input#address-finder"
$ws.Range("C22").Value = "3710 Pio Pico St"
$ws.Range("D22").Value = 1000
$ws.Range("E22").Value = 2000

# Row 23
$ws.Range("A23").Value = "presskey"
$ws.Range("B23").Value = "input#address-finder"
$ws.Range("C23").Value = "Enter"
$ws.Range("D23").Value = 1000
$ws.Range("E23").Value = 2000

# Row 24
$ws.Range("A24").Value = "scroll"
$ws.Range("B24").Value = "Continue button"
$ws.Range("D24").Value = 1000
$ws.Range("E24").Value = 2000

# Row 25
$ws.Range("A25").Value = "click"
$ws.Range("B25").Value = "Phone Number input field"
$ws.Range("D25").Value = 1000
$ws.Range("E25").Value = 2000

# Row 26
$ws.Range("A26").Value = "type"
$ws.Range("B26").Value = "Phone Number input field"
$ws.Range("C26").Value = 6142273098
$ws.Range("D26").Value = 1000
$ws.Range("E26").Value = 2000

# Row 27
$ws.Range("A27").Value = "click"
$ws.Range("B27").Value = "Continue button"
$ws.Range("D27").Value = 1000
# --- Recreate the hyperlinks on their new cells ---
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.bluenile.com/jewelry/necklaces/lab-grown-diamond-cushion-cut-solitaire-pendant-in-14k-white-gold-1-2-ct-tw-f-g-vs2-si1-item-202314") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:mellina@gmail.com") | Out-Null

# --- Re-apply the "Hyperlink" style (matching the original style index) to the two
#     cells that now contain hyperlink text, then clean up the scratch cell. ---
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

Write-Output "Applied bluenile.xlsx update successfully"